$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 -> "Group_1" row (previously row 4's EAN values moved up / re-derived)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Group_1"
$ws.Range("B2").Value = 4987241155729
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = "exclude"
$ws.Range("G2").Value = ""
$ws.Rows.Item(2).RowHeight = 14.95

# ---------------------------------------------------------------------------
# Row 3 -> "Group_2" row
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Group_2"
$ws.Range("B3").Value = "4901872061280, 4901872099122, 4901872963461"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "1,2"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "exclude"
$ws.Range("G3").Value = "Secondary(Side net), Secondary Self Skin Care"
$ws.Rows.Item(3).RowHeight = 41.95

# ---------------------------------------------------------------------------
# Row 4 -> "Group_3" row
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "Group_3"
$ws.Range("B4").Value = "4901872963300, 4901872963461, 4901872049882, 4987241155736"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "5,6,7"
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "exclude"
$ws.Range("G4").Value = "Secondary(Side net), Secondary Self Skin Care"
$ws.Rows.Item(4).RowHeight = 55.45

# ---------------------------------------------------------------------------
# Row 5 -> "Group_4" row (only the "stacking" value changes to "exclude")
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Group_4"
$ws.Range("B5").Value = "4901872049912, 4901872049844, 4901872049790, 4901301346995"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "2,3"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "exclude"
$ws.Range("G5").Value = ""

# ---------------------------------------------------------------------------
# A trailing, otherwise empty row appears at the very bottom of the sheet
# (artifact of the sheet being re-saved against the full row range).
# ---------------------------------------------------------------------------
$ws.Rows.Item(1048576).RowHeight = 12.8

# ---------------------------------------------------------------------------
# Column widths were all nudged slightly wider.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.666666666666668
$ws.Columns.Item(2).ColumnWidth = 22.833333333333332
$ws.Columns.Item(3).ColumnWidth = 22.666666666666668
$ws.Columns.Item(4).ColumnWidth = 27.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.0
$ws.Columns.Item(6).ColumnWidth = 15.5
$ws.Columns.Item(7).ColumnWidth = 42.833333333333336
$ws.Columns.Item(8).ColumnWidth = 9.5

# ---------------------------------------------------------------------------
# Active selection moved to G2.
# ---------------------------------------------------------------------------
$null = $ws.Range("G2").Select()
